# adiciona fonte das kbas e areas do pat
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 22: "Áreas - PAT" -> fully described PAT area, with its source details filled in
$ws.Range("A22").Value2 = "Áreas - PAT (Plano de Ação Territorial Chapada Diamantina-Serra da Jiboia)"
$ws.Range("D22").Value2 = "INEMA/BA"
$ws.Range("E22").Value2 = 2020
$ws.Range("H22").Value2 = "http://www.inema.ba.gov.br/plano-de-acao-territorial-pat-chapada-diamantina-serra-da-jiboia/"

# Row 23: "KBAs" -> fill in its source details
$ws.Range("D23").Value2 = "Key Biodiversity Area Partnership"
$ws.Range("H23").Value2 = "https://www.keybiodiversityareas.org/"
$ws.Range("I23").Value2 = 45013
$ws.Range("I23").NumberFormat = $ws.Range("I16").NumberFormat

# Row 24 ("Labels"/"Texto" placeholder row) is no longer needed, remove it entirely
$ws.Rows.Item(24).Delete()

# Update selection to reflect where editing finished
$ws.Range("D23").Select()
